$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new record is inserted as row 33, pushing the existing records
# (old rows 33-142) down by one row (new rows 34-143).
$ws.Rows(33).Insert()

# After the insert, row 34 holds what used to be row 33's data (the
# insert shifted it down). Copy that record into the newly blank row 33
# so the new row starts as a duplicate of the old row 33 record...
$ws.Range("A34:R34").Copy()
$ws.Range("A33:R33").PasteSpecial()

# ...then overwrite the two fields that differ for the new record.
$ws.Range("D33").Value = 44592
$ws.Range("J33").Value = 100
